# Multiplied DVF errors by pixel size (supposed isotropic) to get the correct values
# -> add a new "pixel_size_mm" column (I) with header + value, matching the
#    formatting already used by the other header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell I1, with the same (bold + bottom-border) style as H1.
$ws.Range("I1").Value = "pixel_size_mm"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats

# New data cell I2: pixel size (mm) used to rescale the DVF errors.
$ws.Range("I2").Value = 1.818

# Widen column I so the new header text is fully visible.
$ws.Columns.Item(9).ColumnWidth = 12.5

# Match the author's final cursor position.
$ws.Range("K3").Select()
